$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "PDB molecule" column (column D) entirely - shifts
# PDB filename / Is model left by one column.
$ws.Range("D1").EntireColumn.Delete()

# Select the column that now holds "PDB filename" (was column E, now D)
# to match the resulting selection state.
$ws.Range("D1:D1048576").Select()
